# BaoCaoTuan.docx formatting pass
# - Apply Times New Roman + bold/size to the title and every "Báo cáo tuần N"
#   heading, body line and bullet item.
# - Merge runs that got split across edits ("Báo cáo tuần X", "Cài git
#   bash...", "Trong tuần này... được:") back into a single run by doing a
#   Find/Replace with the exact same (unsplit) text.
# - Re-stamp the two bare `<w:p/>` paragraphs with paragraph-mark formatting.

$d = $word.ActiveDocument

function Set-RunFormatting($para, [int]$halfPoints, [bool]$bold) {
    $rng = $para.Range
    $rng.Font.Name = "Times New Roman"
    $rng.Font.NameBi = "Times New Roman"
    $rng.Font.Size = $halfPoints / 2
    $rng.Font.SizeBi = $halfPoints / 2
    if ($bold) {
        $rng.Bold = 1
        $rng.BoldBi = 1
    }
}

function Set-EmptyParaFormatting($para, [int]$halfPoints) {
    # Paragraphs with no runs silently drop most Font writes in this host,
    # so stamp a throwaway run, format the (now real) range, then delete
    # just the inserted character -- the paragraph-mark rPr left behind
    # keeps the formatting.
    $start = $para.Range.Start
    $para.Range.InsertBefore("X")
    $p2 = $d.Paragraphs.Item($start + 1)
    $p2 = $d.Range($start, $start + 1).Paragraphs.Item(1)
    $rng = $p2.Range
    $rng.Font.Name = "Times New Roman"
    $rng.Font.NameBi = "Times New Roman"
    $rng.Font.Size = $halfPoints / 2
    $rng.Font.SizeBi = $halfPoints / 2
    $charRange = $d.Range($start, $start + 1)
    $charRange.Text = ""
}

# ---------------------------------------------------------------------
# 1. Normalize runs that were split across multiple <w:r> elements.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Báo cáo tuần 2", $true, $false, $false, $false, $false, $true, 1, $false, "Báo cáo tuần 2", 2) | Out-Null
$d.Content.Find.Execute("Báo cáo tuần 3", $true, $false, $false, $false, $false, $true, 1, $false, "Báo cáo tuần 3", 2) | Out-Null
$d.Content.Find.Execute("Báo cáo tuần 4", $true, $false, $false, $false, $false, $true, 1, $false, "Báo cáo tuần 4", 2) | Out-Null
$d.Content.Find.Execute("Cài git bash, cách sử dụng git.", $true, $false, $false, $false, $false, $true, 1, $false, "Cài git bash, cách sử dụng git.", 2) | Out-Null
$d.Content.Find.Execute("Trong tuần này, em biết thêm được:", $true, $false, $false, $false, $false, $true, 1, $false, "Trong tuần này, em biết thêm được:", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Apply Times New Roman + size (+ bold where applicable) paragraph by
#    paragraph, using the now-stable paragraph indices.
# ---------------------------------------------------------------------

# 1: title - 14pt, bold
Set-RunFormatting $d.Paragraphs.Item(1) 28 $true
# 2: "Báo cáo tuần 1" - 13pt, bold
Set-RunFormatting $d.Paragraphs.Item(2) 26 $true
# 3: "Trong tuần này, em biết thêm được:" - 13pt
Set-RunFormatting $d.Paragraphs.Item(3) 26 $false
# 4: bullet - "Soạn tên đề tài..."
Set-RunFormatting $d.Paragraphs.Item(4) 26 $false
# 5: bullet - "Chọn công nghệ..."
Set-RunFormatting $d.Paragraphs.Item(5) 26 $false
# 6: empty bold paragraph (already has b/bCs) - just needs font/size added
Set-RunFormatting $d.Paragraphs.Item(6) 26 $true
# 7: "Báo cáo tuần 2" - 13pt, bold
Set-RunFormatting $d.Paragraphs.Item(7) 26 $true
# 8: "Trong tuần này, em biết thêm được:"
Set-RunFormatting $d.Paragraphs.Item(8) 26 $false
# 9: "Cài git bash, cách sử dụng git."
Set-RunFormatting $d.Paragraphs.Item(9) 26 $false
# 10: "Tạo repo trong github"
Set-RunFormatting $d.Paragraphs.Item(10) 26 $false
# 11: bare <w:p/> -> gets pPr/rPr only (no bold)
Set-EmptyParaFormatting $d.Paragraphs.Item(11) 26
# 12: "Báo cáo tuần 3" (bookmarkStart lives here too) - 13pt, bold
Set-RunFormatting $d.Paragraphs.Item(12) 26 $true
# 13: "Trong tuần này, em biết thêm được:"
Set-RunFormatting $d.Paragraphs.Item(13) 26 $false
# 14: "Chạy được trang chủ framework Laravel"
Set-RunFormatting $d.Paragraphs.Item(14) 26 $false
# 15: "Chạy được layout home"
Set-RunFormatting $d.Paragraphs.Item(15) 26 $false
# 16: bare <w:p/> -> gets pPr/rPr only (no bold)
Set-EmptyParaFormatting $d.Paragraphs.Item(16) 26
# 17: "Báo cáo tuần 4" - 13pt, bold
Set-RunFormatting $d.Paragraphs.Item(17) 26 $true
# 18: "Trong tuần này, em biết thêm được:"
Set-RunFormatting $d.Paragraphs.Item(18) 26 $false
# 19: "Chạy được giao diện admin"
Set-RunFormatting $d.Paragraphs.Item(19) 26 $false
# 20: "Chạy được trang chủ project "
Set-RunFormatting $d.Paragraphs.Item(20) 26 $false

Write-Host "done"
